# Add 2022-Q3 data
# -----------------
# 1) Update the "总计" (totals) sheet with a new top row for the new quarter,
#    cascading the other quarters' 日期/持有数量/持有市值 figures down by one
#    row and appending the oldest quarter into a brand-new last row.
# 2) Insert a brand-new worksheet "2022-Q3" right after "总计" (before "2022-Q2"),
#    built from the "2022-Q2" sheet's layout/formatting, populated with the new
#    quarter's single fund holding.
# All the other quarterly sheets are left untouched content-wise; they simply
# shift one position to the right because of the newly inserted sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: "总计" sheet - add a new quarter row and cascade the older rows
# down. Note the "A" column here is just a static per-row position index
# (0,1,2,...) that never changes value -- only the 日期/持有数量/持有市值
# columns (B/C/D) cascade down one row at a time, with the oldest quarter
# falling into a brand-new last row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Clone row 7's format+values into a brand-new row 8 (keeps column A's style,
# i.e. bold/bordered index cell, consistent with the rest of the table).
$total.Range("A7:D7").Copy($total.Range("A8"))
$total.Range("A8").Value = 6

# Cascade B/C/D down from row r into row r+1, walking bottom-up so a row's
# old values are read before they get overwritten. This leaves row 2 free
# for the brand-new 2022-Q3 figures.
for ($r = 7; $r -ge 2; $r--) {
    $dstRow = $r + 1
    $total.Cells.Item($dstRow, 2).Value = $total.Cells.Item($r, 2).Value2
    $total.Cells.Item($dstRow, 3).Value = $total.Cells.Item($r, 3).Value2
    $total.Cells.Item($dstRow, 4).Value = $total.Cells.Item($r, 4).Value2
}

# Write the brand-new 2022-Q3 totals into row 2 (A2's running index of 0 is
# already correct and untouched).
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.03

# ---------------------------------------------------------------------------
# Step 2: insert a new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($null, $afterSheet, $null, $null)
$newSheet.Name = "2022-Q3"

# Re-fetch sheets by name (index-based handles go stale once sheets are
# inserted/renamed).
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Item("2022-Q3")

# Clone the header row + one data row from "2022-Q2" so "2022-Q3" starts out
# with identical column headers, number formats (text-like numeric strings)
# and cell styles.
$q2.Range("A1:H2").Copy($q3.Range("A1"))

# Make sure the text-like numeric columns stay text (matches the source
# workbook's convention of storing these figures as strings) before writing
# the real 2022-Q3 values over the templated row.
$q3.Range("B2:G2").NumberFormat = "@"
$q3.Range("B2").Value = "516770"
$q3.Range("C2").Value = "华泰柏瑞中证动漫游戏ETF"
$q3.Range("D2").Value = "0.99"
$q3.Range("E2").Value = "96.39"
$q3.Range("F2").Value = "3.06"
$q3.Range("G2").Value = "0.0303"
$q3.Range("H2").Value = 10

# ---------------------------------------------------------------------------
# Restore the originally-active tab ("2021-Q1"), since inserting/renaming
# sheets above left the newly-added sheet focused.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
